# "przeście na sterowanie z przeglądarki"
# Update the tracked advertiser/display count and move the active selection,
# matching the author's edit to Arkusz1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")
$ws.Activate() | Out-Null

# A2: 93931 -> 88278
$ws.Range("A2").Value = 88278

# Active selection moves from F6 to K23
$ws.Range("K23").Select() | Out-Null
